# "adding new progress as of date 04 nov 2025"
#
# The "Training Dashboard" sheet tracks, for each cargo training row, how
# many days remain until the training expires (column H, "PERIOD TO
# EXPIRE") and when that figure was last computed (column I, "LAST
# UPDATE"). Refreshing the report one day later (03-Nov-2025 ->
# 04-Nov-2025) drops every remaining-day count by one and bumps the
# "LAST UPDATE" stamp to the new date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# row -> new "PERIOD TO EXPIRE" value (each one day less than before)
$periodToExpire = @{
    3 = 182
    4 = 190
    5 = 191
    6 = 195
    7 = 170
    8 = 189
    9 = 174
}

$newLastUpdate = "04-Nov-2025"

foreach ($row in $periodToExpire.Keys) {
    $ws.Cells.Item($row, 8).Value = $periodToExpire[$row]

    # Leading apostrophe forces the Excel text-entry semantics so the
    # date-looking string is stored as literal text (matching the
    # existing "LAST UPDATE" column, which holds plain date strings,
    # not real date serials) instead of being parsed into a date value.
    $ws.Cells.Item($row, 9).Value = "'" + $newLastUpdate
}
